# Update cryptocurrency price/volume data to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/URL/percentage cells - safe to assign directly.
$plainUpdates = @(
    @{ Addr = 'D2'; Val = '62.809.84' }
    @{ Addr = 'E2'; Val = '  +0.94%  ' }
    @{ Addr = 'D3'; Val = '3.433.69' }
    @{ Addr = 'E3'; Val = '  +0.81%  ' }
    @{ Addr = 'E4'; Val = '  -0.16%  ' }
    @{ Addr = 'E5'; Val = '  +0.02%  ' }
    @{ Addr = 'E6'; Val = '  -0.44%  ' }
    @{ Addr = 'E7'; Val = '  -1.88%  ' }
    @{ Addr = 'E8'; Val = '  -0.18%  ' }
    @{ Addr = 'E9'; Val = '  +2.20%  ' }
    @{ Addr = 'E10'; Val = '  +9.43%  ' }
    @{ Addr = 'E11'; Val = '  -0.56%  ' }
    @{ Addr = 'E12'; Val = '  -0.18%  ' }
    @{ Addr = 'E13'; Val = '  -1.68%  ' }
    @{ Addr = 'E14'; Val = '  +0.21%  ' }
    @{ Addr = 'D15'; Val = '3.425.14' }
    @{ Addr = 'E15'; Val = '  +0.33%  ' }
    @{ Addr = 'D16'; Val = '62.547.07' }
    @{ Addr = 'E16'; Val = '  +0.76%  ' }
    @{ Addr = 'E17'; Val = '  +0.44%  ' }
    @{ Addr = 'E18'; Val = '  -1.65%  ' }
    @{ Addr = 'E19'; Val = '  +16.45%  ' }
    @{ Addr = 'E20'; Val = '  -3.12%  ' }
    @{ Addr = 'E21'; Val = '  +1.82%  ' }
    @{ Addr = 'E22'; Val = '  +1.80%  ' }
    @{ Addr = 'E23'; Val = '  -2.96%  ' }
    @{ Addr = 'E24'; Val = '  +0.13%  ' }
    @{ Addr = 'E25'; Val = '  +1.32%  ' }
    @{ Addr = 'E26'; Val = '  -0.11%  ' }
    @{ Addr = 'E27'; Val = '  -5.49%  ' }
    @{ Addr = 'E28'; Val = '  +3.98%  ' }
    @{ Addr = 'B29'; Val = 'InjectiveProtocol' }
    @{ Addr = 'C29'; Val = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Addr = 'E29'; Val = '  +6.47%  ' }
    @{ Addr = 'B30'; Val = 'Toncoin' }
    @{ Addr = 'C30'; Val = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Addr = 'E30'; Val = '  +4.21%  ' }
    @{ Addr = 'E31'; Val = '  -0.99%  ' }
    @{ Addr = 'E32'; Val = '  -1.53%  ' }
    @{ Addr = 'E33'; Val = '  -3.45%  ' }
    @{ Addr = 'E34'; Val = '  +0.05%  ' }
    @{ Addr = 'E35'; Val = '  -1.17%  ' }
    @{ Addr = 'E36'; Val = '  -1.07%  ' }
    @{ Addr = 'E37'; Val = '  +0.05%  ' }
    @{ Addr = 'B38'; Val = 'Stacks' }
    @{ Addr = 'C38'; Val = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Addr = 'E38'; Val = '  -0.33%  ' }
    @{ Addr = 'B39'; Val = 'TheGraph' }
    @{ Addr = 'C39'; Val = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' }
    @{ Addr = 'E39'; Val = '  +13.16%  ' }
    @{ Addr = 'E40'; Val = '  -2.65%  ' }
    @{ Addr = 'E41'; Val = '  +3.29%  ' }
    @{ Addr = 'E42'; Val = '  -0.05%  ' }
    @{ Addr = 'E43'; Val = '  -2.33%  ' }
    @{ Addr = 'E44'; Val = '  -1.33%  ' }
    @{ Addr = 'E45'; Val = '  -1.45%  ' }
    @{ Addr = 'E46'; Val = '  -0.61%  ' }
    @{ Addr = 'E47'; Val = '  -2.64%  ' }
    @{ Addr = 'D48'; Val = '2.106.80' }
    @{ Addr = 'E48'; Val = '  -2.24%  ' }
    @{ Addr = 'B49'; Val = 'ThetaToken' }
    @{ Addr = 'C49'; Val = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta' }
    @{ Addr = 'E49'; Val = '  +3.03%  ' }
    @{ Addr = 'B50'; Val = 'ApeXProtocol' }
    @{ Addr = 'C50'; Val = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ Addr = 'E50'; Val = '  -1.37%  ' }
    @{ Addr = 'E51'; Val = '  +29.04%  ' }
)
foreach ($u in $plainUpdates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# Price cells that look numeric (e.g. "0.999", "45.28") must stay as text
# to match the source sheet, which stores them as inline strings. Force the
# "Text" number format before assigning, then restore the default cell style
# so no stray style index gets attached to the cell.
$textUpdates = @(
    @{ Addr = 'D4'; Val = '0.999' }
    @{ Addr = 'D5'; Val = '407.01' }
    @{ Addr = 'D6'; Val = '130.43' }
    @{ Addr = 'D7'; Val = '0.596' }
    @{ Addr = 'D8'; Val = '0.998' }
    @{ Addr = 'D9'; Val = '0.693' }
    @{ Addr = 'D10'; Val = '0.139' }
    @{ Addr = 'D11'; Val = '42.02' }
    @{ Addr = 'D13'; Val = '8.42' }
    @{ Addr = 'D14'; Val = '19.83' }
    @{ Addr = 'D17'; Val = '11.58' }
    @{ Addr = 'D18'; Val = '1.02' }
    @{ Addr = 'D19'; Val = '0.0000156' }
    @{ Addr = 'D20'; Val = '3.18' }
    @{ Addr = 'D21'; Val = '84.37' }
    @{ Addr = 'D22'; Val = '313.56' }
    @{ Addr = 'D23'; Val = '12.81' }
    @{ Addr = 'D24'; Val = '3.18' }
    @{ Addr = 'D25'; Val = '4.75' }
    @{ Addr = 'D26'; Val = '29.76' }
    @{ Addr = 'D27'; Val = '8.12' }
    @{ Addr = 'D28'; Val = '7.79' }
    @{ Addr = 'D29'; Val = '45.28' }
    @{ Addr = 'D30'; Val = '2.74' }
    @{ Addr = 'D33'; Val = '11.38' }
    @{ Addr = 'D36'; Val = '51.89' }
    @{ Addr = 'D37'; Val = '0.998' }
    @{ Addr = 'D38'; Val = '2.98' }
    @{ Addr = 'D39'; Val = '0.323' }
    @{ Addr = 'D40'; Val = '3.34' }
    @{ Addr = 'D41'; Val = '141.89' }
    @{ Addr = 'D42'; Val = '0.126' }
    @{ Addr = 'D45'; Val = '16.84' }
    @{ Addr = 'D47'; Val = '21.27' }
    @{ Addr = 'D49'; Val = '1.98' }
    @{ Addr = 'D50'; Val = '2.32' }
)
foreach ($u in $textUpdates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    $cell.Style = "Normal"
}
